$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.153.02'
$ws.Range("E2").Value = '  +4.50%  '
$ws.Range("D3").Value = '2.218.81'
$ws.Range("E3").Value = '  +1.86%  '
$ws.Range("E4").Value = '  +0.01%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '260.05'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +2.97%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '82.80'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +12.68%  '
$ws.Range("E7").Value = '  +4.29%  '
$ws.Range("E8").Value = '  -0.08%  '
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.607'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +4.29%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '44.00'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +7.87%  '
$ws.Range("E11").Value = '  +2.65%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.06'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +4.47%  '
$ws.Range("E13").Value = '  +2.75%  '
$ws.Range("D14").Value = '2.553.49'
$ws.Range("E14").Value = '  +1.96%  '
$ws.Range("E15").Value = '  +2.68%  '
$ws.Range("D16").Value = '2.220.94'
$ws.Range("E16").Value = '  +1.65%  '
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.778'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +1.73%  '
$ws.Range("D18").Value = '44.060.24'
$ws.Range("E19").Value = '  +1.84%  '
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '71.28'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +0.92%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.01'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +2.47%  '
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.38'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +10.14%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '233.03'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +2.80%  '
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.35'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -1.40%  '
$ws.Range("E26").Value = '  +3.20%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '41.49'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +12.71%  '
$ws.Range("E29").Value = '  +2.98%  '
$ws.Range("E30").Value = '  +0.12%  '
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '172.84'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +2.16%  '
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '20.66'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +3.26%  '
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0878'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +9.10%  '
$ws.Range("E34").Value = '  +4.35%  '
$ws.Range("E35").Value = '  +7.97%  '
$ws.Range("E36").Value = '  +2.37%  '
$ws.Range("E37").Value = '  +9.03%  '
$ws.Range("E38").Value = '  +6.52%  '
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '13.60'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +14.00%  '
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.95'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +19.70%  '
$ws.Range("E41").Value = '  +3.46%  '
$ws.Range("E42").Value = '  +8.97%  '
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '63.37'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +7.18%  '
$ws.Range("E44").Value = '  +3.43%  '
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '103.16'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0988'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +2.05%  '
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.33'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +0.48%  '
$ws.Range("E48").Value = '  +29.78%  '
$ws.Range("E49").Value = '  +3.45%  '
$ws.Range("E50").Value = '  +4.38%  '
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.442'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -5.51%  '
